$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.624.53'
$ws.Range('E2').Value = '  -1.87%  '
$ws.Range('D3').Value = '2.463.39'
$ws.Range('E3').Value = '  -2.59%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.80'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.46%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.551'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.11%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -3.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.88'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0782'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.55%  '
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('E13').Value = '  -4.83%  '
$ws.Range('D14').Value = '2.844.91'
$ws.Range('E14').Value = '  -2.52%  '
$ws.Range('D15').Value = '2.450.04'
$ws.Range('E15').Value = '  -4.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.64'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -6.67%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.787'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.92%  '
$ws.Range('D18').Value = '41.591.84'
$ws.Range('E18').Value = '  -1.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.37'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -6.26%  '
$ws.Range('E20').Value = '  -3.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.55'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.77'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.35'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.06%  '
$ws.Range('E24').Value = '  -3.99%  '
$ws.Range('E25').Value = '  -5.64%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.67'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.75%  '
$ws.Range('E28').Value = '  -4.77%  '
$ws.Range('E29').Value = '  -3.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.33'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '154.12'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.61'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.62'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('E34').Value = '  -7.91%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0756'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.82%  '
$ws.Range('E36').Value = '  -4.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.37'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.18%  '
$ws.Range('E38').Value = '  -6.79%  '
$ws.Range('E39').Value = '  -5.06%  '
$ws.Range('E40').Value = '  -3.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.25'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.52%  '
$ws.Range('E43').Value = '  -0.06%  '
$ws.Range('D44').Value = '1.986.30'
$ws.Range('E44').Value = '  +1.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0285'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.91%  '
$ws.Range('E46').Value = '  -6.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.68'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.43%  '
$ws.Range('D48').Value = '2.702.31'
$ws.Range('E48').Value = '  -2.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '76.50'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.48%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.182'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.25%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '97.72'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.57%  '
